# Apply backtest-results update:
#  - Ticker renamed: TSLA -> NFE, AAPL -> PLUG (column A, rows 2-11)
#  - Final Balance (column E) and Cumulative Returns % (column F) recomputed
#    for every data row (2-11) following the LSTM pipeline fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename tickers -------------------------------------------------------
for ($r = 2; $r -le 11; $r++) {
    $ticker = $ws.Cells.Item($r, 1).Value2
    if ($ticker -eq "TSLA") {
        $ws.Cells.Item($r, 1).Value2 = "NFE"
    } elseif ($ticker -eq "AAPL") {
        $ws.Cells.Item($r, 1).Value2 = "PLUG"
    }
}

# --- Update Final Balance (E) / Cumulative Returns % (F) ------------------
# row => [Final Balance, Cumulative Returns %]
$updates = @{
    2  = @(-17.14741545364173, -101.7147415453642)
    3  = @(-112.538466017716,  -111.2538466017716)
    4  = @(1000,                0)
    5  = @(49.39116001017241,  -95.06088399898276)
    6  = @(610.707381687405,   -38.9292618312595)
    7  = @(414.7290100177484,  -58.52709899822516)
    8  = @(685.6608408487546,  -31.43391591512454)
    9  = @(1000,                0)
    10 = @(391.4648580423796,  -60.85351419576204)
    11 = @(385.2095983192128,  -61.47904016807873)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 5).Value2 = $vals[0]
    $ws.Cells.Item($row, 6).Value2 = $vals[1]
}
